# The underlying data rows (2-18) were reordered (the export was re-sorted),
# while every row's own full content (all columns A:AY) stayed exactly the
# same. Move whole rows with Range.Copy (instead of reading/writing typed
# values) so text that merely looks numeric or date-like (e.g. "20",
# "2020-01-10") keeps its original inlineStr/string representation instead
# of being reinterpreted by the host as a number or date.
#
# Because several destination rows are also source rows for other
# destinations, first stage every source row far below the used range,
# then copy from the staging rows into their final destinations, then
# clear the staging rows again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol = "AY"
$stagingOffset = 200

# Mapping: new row number -> old row number (where that row's data used to live).
$mapping = [ordered]@{
    2  = 8
    3  = 2
    4  = 5
    5  = 7
    6  = 10
    7  = 11
    8  = 12
    9  = 14
    10 = 15
    11 = 16
    12 = 17
    13 = 18
    14 = 3
    15 = 4
    16 = 6
    17 = 9
    18 = 13
}

# Every old row that feeds some new row, without duplicates.
$oldRows = $mapping.Values | Sort-Object -Unique

# Stage each distinct source row well away from the live data first.
foreach ($oldRow in $oldRows) {
    $staged = $oldRow + $stagingOffset
    $src = $ws.Range($firstCol + $oldRow + ":" + $lastCol + $oldRow)
    $dst = $ws.Range($firstCol + $staged + ":" + $lastCol + $staged)
    $src.Copy($dst)
}

# Now copy from the staged rows into their real destinations. Clear the
# destination first: Copy only overwrites cells that the source actually
# has a value for, so a truly-blank source cell would otherwise leave
# behind whatever was previously sitting in the destination row.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $staged = $oldRow + $stagingOffset
    $src = $ws.Range($firstCol + $staged + ":" + $lastCol + $staged)
    $dst = $ws.Range($firstCol + $newRow + ":" + $lastCol + $newRow)
    $dst.Clear()
    $src.Copy($dst)
}

# Clean up the staging rows so they don't linger as extra used range.
foreach ($oldRow in $oldRows) {
    $staged = $oldRow + $stagingOffset
    $stagedRange = $ws.Range($firstCol + $staged + ":" + $lastCol + $staged)
    $stagedRange.Clear()
}
